$d = $word.ActiveDocument

# Remove the "Last updated by Randy Peterson" paragraph (style "Author"),
# which sits right above the "Last updated on ..." paragraph (style "Date").
$d.Paragraphs(2).Range.Delete()

# Update the day-of-month in the remaining date line from 2 to 3
# (the "Last updated on October 2, 2018" paragraph, now paragraph 2).
$dateRng = $d.Paragraphs(2).Range.Duplicate
$dateRng.Find.Execute("2,", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "3,", 2)
